$wb = $excel.ActiveWorkbook

# Duplicate the most recent weekly ranking sheet ("2025-12-01") to use as the
# template for the new week, placing the copy immediately after it.
$sourceSheet = $wb.Worksheets.Item("2025-12-01")
$sourceSheet.Copy($null, $sourceSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "2025-12-08"

# rank (column A) stays 1..50 unchanged (copied from the template sheet).
# Update title / author / latest_episode for every ranked entry.

$newSheet.Range("B2").Value = "新米オッサン冒険者、最強パーティに死ぬほど鍛えられて無敵になる"
$newSheet.Range("C2").Value = "漫画：荻野ケン 原作：岸馬きらく キャラクター原案：Tea"
$newSheet.Range("D2").Value = "第73話 後編"

$newSheet.Range("B3").Value = "転生コロシアム～最弱スキルで最強の女たちを攻略して奴隷ハーレム作ります～"
$newSheet.Range("C3").Value = "zunta(作画) はらわたさいぞう(原作)"
$newSheet.Range("D3").Value = "第34話：プロのテク①"

$newSheet.Range("B4").Value = "生徒会にも穴はある！"
$newSheet.Range("C4").Value = "むちまろ"
$newSheet.Range("D4").Value = "第141話	虎丸がんばる！"

$newSheet.Range("B5").Value = "とんでもスキルで異世界放浪メシ"
$newSheet.Range("C5").Value = "赤岸K（漫画） 江口連（原作） 雅（キャラクター原案）"
$newSheet.Range("D5").Value = "第57話　「誘惑だらけの朝市」"

$newSheet.Range("B6").Value = "時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―"
$newSheet.Range("C6").Value = "光永康則"
$newSheet.Range("D6").Value = "第７２話『先端停止』①"

$newSheet.Range("B7").Value = "地元のいじめっ子達に仕返ししようとしたら、別の戦いが始まった。"
$newSheet.Range("C7").Value = "マツモトケンゴ"
$newSheet.Range("D7").Value = "第７０話　鈍感な戦いが始まった"

$newSheet.Range("B8").Value = "王子様の友達"
$newSheet.Range("C8").Value = "すけろく(著者)"
$newSheet.Range("D8").Value = "第32話"

$newSheet.Range("B9").Value = "帰ってください！ 阿久津さん"
$newSheet.Range("C9").Value = "長岡太一(著者)"
$newSheet.Range("D9").Value = "第199話"

$newSheet.Range("B10").Value = "元・世界１位のサブキャラ育成日記 ～廃プレイヤー、異世界を攻略中！～"
$newSheet.Range("C10").Value = "沢村治太郎(原作) 前田理想(漫画) まろ(キャラクター原案)"
$newSheet.Range("D10").Value = "第81話その1"

$newSheet.Range("B11").Value = "実は俺、最強でした？"
$newSheet.Range("C11").Value = "原作：澄守 彩 漫画：高橋 愛"
$newSheet.Range("D11").Value = "第133話　秘密"

$newSheet.Range("B12").Value = "黒幕一家に転生したけど原作無視して独立する"
$newSheet.Range("C12").Value = "空野進 赤村晃人 笠間三四郎 るろお"
$newSheet.Range("D12").Value = "第５話　ゲス子爵を成敗して独立する（３）"

$newSheet.Range("B13").Value = "いとこのこ"
$newSheet.Range("C13").Value = "いぬちく(著者)"
$newSheet.Range("D13").Value = "休載イラスト"

$newSheet.Range("B14").Value = "怠惰な悪辱貴族に転生した俺、シナリオをぶっ壊したら規格外の魔力で最凶になった"
$newSheet.Range("C14").Value = "菊池快晴(原作) 小田童馬(作画) 桑島黎音(キャラクター原案)"
$newSheet.Range("D14").Value = "第17話前半"

$newSheet.Range("B15").Value = "貞操逆転世界で頼めばヤれると噂の俺"
$newSheet.Range("C15").Value = "澄田佑貴(漫画) aaa168（スリーエー）(原作)"
$newSheet.Range("D15").Value = "第2話"

$newSheet.Range("B16").Value = "田舎で恋は難しい!?"
$newSheet.Range("C16").Value = "ねこうめ(著者)"
$newSheet.Range("D16").Value = "第2話"

$newSheet.Range("B17").Value = "異世界魔王と召喚少女の奴隷魔術"
$newSheet.Range("C17").Value = "原作：むらさきゆきや 漫画：福田直叶 キャラクター原案：鶴崎貴大"
$newSheet.Range("D17").Value = "第131話　幕間（前編）"

$newSheet.Range("B18").Value = "女友達は頼めば意外とヤらせてくれる"
$newSheet.Range("C18").Value = "ろくろ(漫画) 鏡遊(原作)"
$newSheet.Range("D18").Value = "第27話①"

$newSheet.Range("B19").Value = "勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～"
$newSheet.Range("C19").Value = "漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり"
$newSheet.Range("D19").Value = "第５５話　封印を解く器用貧乏（２）"

$newSheet.Range("B20").Value = "このヒーラー、めんどくさい"
$newSheet.Range("C20").Value = "丹念に発酵(著者)"
$newSheet.Range("D20").Value = "第92話：特訓"

$newSheet.Range("B21").Value = "Ｓ級ギルドを追放されたけど、実は俺だけドラゴンの言葉がわかるので、気付いたときには竜騎士の頂点を極めてました。"
$newSheet.Range("C21").Value = "ひそな(漫画) 三木なずな(原作) 白狼(キャラクター原案)"
$newSheet.Range("D21").Value = "第41話-2"

$newSheet.Range("B22").Value = "リビルドワールド"
$newSheet.Range("C22").Value = "綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)"
$newSheet.Range("D22").Value = "第75話④"

$newSheet.Range("B23").Value = "バキ外伝 烈海王は異世界転生しても一向にかまわんッッ"
$newSheet.Range("C23").Value = "板垣恵介 猪原賽 陸井栄史"
$newSheet.Range("D23").Value = "第85話　不入虎穴(虎穴に入らずんば虎子を得ず)"

$newSheet.Range("B24").Value = "貞操逆転世界の童貞辺境領主騎士"
$newSheet.Range("C24").Value = "柳瀬こたつ（漫画） 道造（原作） めろん２２（キャラクター原案）"
$newSheet.Range("D24").Value = "第11話　やむなき犠牲（前編）"

$newSheet.Range("B25").Value = "異世界食堂　洋食のねこや"
$newSheet.Range("C25").Value = "犬塚惇平(ヒーロー文庫／イマジカインフォス)(原作) ヤミザワ(漫画) モロザワ(漫画) エナミカツミ(キャラクター原案)"
$newSheet.Range("D25").Value = "第42話➁"

$newSheet.Range("B26").Value = "世界最強の魔女、始めました 〜私だけ『攻略サイト』を見れる世界で自由に生きます〜"
$newSheet.Range("C26").Value = "戸賀 環 坂木持丸 riritto"
$newSheet.Range("D26").Value = "第57話①　スローライフをしてみた"

$newSheet.Range("B27").Value = "濁る瞳で何を願う ハイセルク戦記"
$newSheet.Range("C27").Value = "トルトネン 創-taro 斎藤八呑"
$newSheet.Range("D27").Value = "第35話 黒の奔流"

$newSheet.Range("B28").Value = "聖者無双"
$newSheet.Range("C28").Value = "漫画：秋風緋色 原作：ブロッコリーライオン キャラクター原案：sime"
$newSheet.Range("D28").Value = "第95話　奴隷の扱い・戦闘準備（前半）"

$newSheet.Range("B29").Value = "小林さんちのメイドラゴン"
$newSheet.Range("C29").Value = "クール教信者"
$newSheet.Range("D29").Value = "第155話"

$newSheet.Range("B30").Value = "よくわからないけれど異世界に転生していたようです"
$newSheet.Range("C30").Value = "内々けやき あし カオミン"
$newSheet.Range("D30").Value = "第144話 よくわからないけれど港町に着いたみたいです（１）"

$newSheet.Range("B31").Value = "異世界メイドの三ツ星グルメ ～現代ごはん作ったら王宮で大バズリしました～"
$newSheet.Range("C31").Value = "モリタ Ｕ４ nima"
$newSheet.Range("D31").Value = "第14話（４）　春とおぼっちゃまとピクニックランチ（４）"

$newSheet.Range("B32").Value = "願ってもない追放後からのスローライフ？ 〜引退したはずが成り行きで美少女ギャルの師匠になったらなぜかめちゃくちゃ懐かれた〜"
$newSheet.Range("C32").Value = "ヤミーゴ(漫画) シュガースプーン。（GA文庫/SBクリエイティブ）(原作) なたーしゃ(キャラクター原案)"
$newSheet.Range("D32").Value = "第9話-1：黄昏の茶会"

$newSheet.Range("B33").Value = "ぽんドロイド！ はまさん"
$newSheet.Range("C33").Value = "はれやまはれぞう(著者)"
$newSheet.Range("D33").Value = "第13話"

$newSheet.Range("B34").Value = "アイドル辞めるけど結婚してくれますか!?"
$newSheet.Range("C34").Value = "三吉汐美(著者)"
$newSheet.Range("D34").Value = "第19話前半"

$newSheet.Range("B35").Value = "治癒魔法の間違った使い方 ~戦場を駆ける回復要員~"
$newSheet.Range("C35").Value = "九我山レキ(漫画) くろかた(原作) ＫｅＧ(キャラクター原案)"
$newSheet.Range("D35").Value = "第84話その2"

$newSheet.Range("B36").Value = "落ちこぼれだった兄が実は最強 ～史上最強の勇者は転生し、学園で無自覚に無双する～"
$newSheet.Range("C36").Value = "村上よしゆき 茨木野 あるてら"
$newSheet.Range("D36").Value = "第４３話　勇者、合体した六邪神将を撃破し、めでたしめでたし（４）"

$newSheet.Range("B37").Value = "衛宮さんちの今日のごはん"
$newSheet.Range("C37").Value = "TAa(漫画) 只野まこと(料理監修) ＴＹＰＥ－ＭＯＯＮ(原作)"
$newSheet.Range("D37").Value = "第77話"

$newSheet.Range("B38").Value = "追放されたチート付与魔術師は 気ままなセカンドライフを謳歌する。"
$newSheet.Range("C38").Value = "六志麻あさ 業務用餅 kisui"
$newSheet.Range("D38").Value = "第７６話"

$newSheet.Range("B39").Value = "くらいあの子としたいこと"
$newSheet.Range("C39").Value = "碇マナツ(著者)"
$newSheet.Range("D39").Value = "第87話"

$newSheet.Range("B40").Value = "まったく最近の探偵ときたら"
$newSheet.Range("C40").Value = "五十嵐正邦(著者)"
$newSheet.Range("D40").Value = "第117話"

$newSheet.Range("B41").Value = "バキ外伝　ガイアとシコルスキー　～ときどきノムラ 二人だけど三人暮らし～"
$newSheet.Range("C41").Value = "板垣恵介 林たかあき"
$newSheet.Range("D41").Value = "第59話 流水の獲物"

$newSheet.Range("B42").Value = "理想のヒモ生活"
$newSheet.Range("C42").Value = "日月ネコ(漫画) 渡辺恒彦（ヒーロー文庫／イマジカインフォス）(原作) 文倉十(キャラクター原案)"
$newSheet.Range("D42").Value = "第90話　その1"

$newSheet.Range("B43").Value = "姫様“拷問”の時間です"
$newSheet.Range("C43").Value = "原作:春原ロビンソン　漫画:ひらけい"
$newSheet.Range("D43").Value = "拷問160"

$newSheet.Range("B44").Value = "経験値貯蓄でのんびり傷心旅行 ～勇者と恋人に追放された戦士の無自覚ざまぁ～"
$newSheet.Range("C44").Value = "奏ヨシキ(著者) 徳川レモン(原作) riritto(キャラクターデザイン)"
$newSheet.Range("D44").Value = "第41話-2"

$newSheet.Range("B45").Value = "ふかふかダンジョン攻略記～俺の異世界転生冒険譚～"
$newSheet.Range("C45").Value = "KAKERU"
$newSheet.Range("D45").Value = "第71話「ファントム・アレイ」（前半）"

$newSheet.Range("B46").Value = "十年目、帰還を諦めた転移者はいまさら主人公になる"
$newSheet.Range("C46").Value = "原作：氷純（「十年目、帰還を諦めた転移者はいまさら主人公になる」MFブックス刊） 漫画：しゅーかま キャラクター原案：あんべよしろう"
$newSheet.Range("D46").Value = "第21話①"

$newSheet.Range("B47").Value = "彼女にしたい女子一位、の隣で見つけたあまりちゃん"
$newSheet.Range("C47").Value = "寝巻ネルゾ(漫画) 裕時悠示(原作) たん旦(キャラクター原案)"
$newSheet.Range("D47").Value = "第7話①「伝えたいこと」"

$newSheet.Range("B48").Value = "陰キャの俺が席替えでS級美少女に囲まれたら秘密の関係が始まった。"
$newSheet.Range("C48").Value = "星野 星野(原作) バラマツヒトミ(漫画) 黒兎 ゆう(キャラクター原案)"
$newSheet.Range("D48").Value = "第8話"

$newSheet.Range("B49").Value = "婚約者に裏切られた錬金術師は、独立して『ざまぁ』します　コミック版"
$newSheet.Range("C49").Value = "漫画/すたひろ 原作/Y.A"
$newSheet.Range("D49").Value = "chapter75【39話②】"

$newSheet.Range("B50").Value = "黄金の経験値"
$newSheet.Range("C50").Value = "原純(原作) 霜月汐(作画) fixro2n(キャラクター原案)"
$newSheet.Range("D50").Value = "第20話（後編）"

$newSheet.Range("B51").Value = "最弱貴族に転生したので悪役たちを集めてみた"
$newSheet.Range("C51").Value = "空野進 sorani ファルまろ"
$newSheet.Range("D51").Value = "第14話　最弱貴族、悪役令嬢を脱がす（４）"

$wb.Worksheets.Item(1).Select()
